$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at row 3 (pushes current rows 3-21 down to 4-22),
#    then populate it with the new "이엔셀" entry that now appears right
#    after "교보스팩16호".
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).Value = "이엔셀"
$ws.Cells.Item(3, 2).Value = "2024.07.23~07.29"
$ws.Cells.Item(3, 3).Value = "13,600~15,300"
$ws.Cells.Item(3, 4).Value = "-"
$ws.Cells.Item(3, 5).Value = 21308
$ws.Cells.Item(3, 6).Value = "NH투자증권"

# 2. The old "이엔셀" row (originally row 15, now shifted to row 16 by the
#    insert above) is removed, since its data has been superseded by the
#    new entry inserted at row 3.
$ws.Rows.Item(16).Delete()
